$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "CB"
$ws.Range("A6").Value = "CB"
$ws.Range("A8").Value = "TS"
$ws.Range("A10").Value = "TS"

$ws.Range("B16").Select()
